$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st tab) - update column F values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3536
$wsExhibit.Range("F4").Value = 146
$wsExhibit.Range("F5").Value = 7026
$wsExhibit.Range("F6").Value = 3296
$wsExhibit.Range("F7").Value = 58
$wsExhibit.Range("F8").Value = 140
$wsExhibit.Range("F13").Value = 18
$wsExhibit.Range("F15").Value = 590
$wsExhibit.Range("F16").Value = 36

# Sheet "全部类型" (4th tab) - update column F values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3536
$wsAll.Range("F5").Value = 146
$wsAll.Range("F6").Value = 7026
$wsAll.Range("F7").Value = 3296
$wsAll.Range("F8").Value = 58
$wsAll.Range("F9").Value = 140
$wsAll.Range("F14").Value = 18
$wsAll.Range("F16").Value = 590
$wsAll.Range("F17").Value = 36
